$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176 (shifts existing rows 176:260 down to 177:261,
# carrying the date-format style of column D down with them).
$ws.Rows(176).Insert()

# Populate the newly inserted row 176 with the new weekly price-report record.
$ws.Range("A176").Value = 5
$ws.Range("B176").Value = "Macroferia Regional de Talca"
$ws.Range("C176").Value = "Maule"
$ws.Range("D176").Value = 45097
$ws.Range("E176").Value = 7
$ws.Range("F176").Value = 100112031
$ws.Range("G176").Value = "Poroto verde"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 150
$ws.Range("K176").Value = 23000
$ws.Range("L176").Value = 23000
$ws.Range("M176").Value = 23000
$ws.Range("N176").Value = "$/malla 25 kilos"
$ws.Range("O176").Value = "Perú"
$ws.Range("P176").Value = 920
$ws.Range("Q176").Value = 25
$ws.Range("R176").Value = "Hortaliza"
